$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:E2").Copy()
$ws.Range("A3").PasteSpecial()
